$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correcao no orcamento de Telecom: valor mensal previsto de R$ 32,000 para R$ 32,500
$ws.Range("B6").Value = "R$ 32,500"

# Atualiza a celula selecionada na planilha
$ws.Range("D3").Select()
